# Update the utilization report with the newly added 64-bit ("only_integer64")
# design numbers. Only column F (DSP) and the row-2 resource-utilization
# figures for LUT (B), FF (D), BRAM (E) and DSP (F) change; LUTRAM (C) and
# the remaining columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.518796920776367
$ws.Range("D2").Value = 14.88063907623291
$ws.Range("E2").Value = 57.85714340209961
$ws.Range("F2").Value = 57.272727966308594

# Column F needs to widen to match the other "wide" columns (same display
# width as columns B, D and E) now that it holds a bigger number.
$ws.Columns.Item(6).ColumnWidth = 10.833333333333334
